$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Description"), shifting the Description column
# (and the whole second/right-hand table) one column to the right.
$ws.Columns("D").Insert()

# New column D is the "Data (2Byte)" column for the left (Host -> FPGA) table.
$ws.Range("D1").Value = "Data (2Byte)"
$ws.Range("D2").Value = "-"

# Two new instruction rows describing GRAM / IRAM writes.
$ws.Range("B3").Value = "0x01"
$ws.Range("C3").Value = "Addr"
$ws.Range("D3").Value = "Data"
$ws.Range("E3").Value = "GRAM Write"

$ws.Range("B4").Value = "0x02"
$ws.Range("C4").Value = "Addr"
$ws.Range("D4").Value = "Data"
$ws.Range("E4").Value = "IRAM Write"

# Adjust column widths to the new layout (values expressed as Excel's
# ColumnWidth, which is ~ OOXML character width minus 5/6).
$ws.Columns("D").ColumnWidth = 17.833333333333332
$ws.Columns("E").ColumnWidth = 25.166666666666668
$ws.Columns("I").ColumnWidth = 18.0
$ws.Columns("J").ColumnWidth = 25.333333333333332

# Move the selection cursor to where the author left it.
[void]$ws.Range("D5").Select()
